$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 53: fill in A/B (previously blank) to match the pattern of surrounding rows
$ws.Range("A53").Value = "7/8/2025(Remote)"
$ws.Range("B53").Value = "Car Tracking Project"

# Row 54: fill in A/B, fix the typo in C (new shared string created here first)
$ws.Range("A54").Value = "7/8/2025(Remote)"
$ws.Range("B54").Value = "Car Tracking Project"
$ws.Range("C54").Value = "secret env variables"

# Row 51: extend the "Remember to change..." text, turn on wrap text, taller row
# (new shared string created here second, so it sorts after "secret env variables")
$ws.Range("C51").Value = "Remember to change the extraction and limit rate to the normal AND remove any emojis from `nyour code base"
$ws.Range("C51").WrapText = $true
$ws.Rows.Item(51).RowHeight = 28.8

# Row 55: fill in A/B, add F = DONE
$ws.Range("A55").Value = "7/8/2025(Remote)"
$ws.Range("B55").Value = "Car Tracking Project"
$ws.Range("F55").Value = "DONE"

# Row 49: add F = DONE
$ws.Range("F49").Value = "DONE"

# Update sheet view (scroll / selection position)
$ws.Application.ActiveWindow.ScrollRow = 42
$ws.Range("C50").Select()
